$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 34 (pushes existing rows 34-46 down to 36-48)
$ws.Rows.Item(34).Resize(2).Insert()

# New row 34: Caqui / Fuyu / Primera
$ws.Range("A34").Value = 10
$ws.Range("B34").Value = "Vega Modelo de Temuco"
$ws.Range("C34").Value = "La Araucanía"
$ws.Range("D34").Value = 44719
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100107
$ws.Range("H34").Value = "Otros"
$ws.Range("I34").Value = 100107001
$ws.Range("J34").Value = "Caqui"
$ws.Range("K34").Value = "Fuyu"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 35
$ws.Range("N34").Value = 20000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 20000
$ws.Range("Q34").Value = "$/bandeja 15 kilos granel"
$ws.Range("R34").Value = "Región de O'Higgins"
$ws.Range("S34").Value = 1333
$ws.Range("T34").Value = 15

# New row 35: Caqui / Mankaki / Primera
$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44719
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100107
$ws.Range("H35").Value = "Otros"
$ws.Range("I35").Value = 100107001
$ws.Range("J35").Value = "Caqui"
$ws.Range("K35").Value = "Mankaki"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 35
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 20000
$ws.Range("P35").Value = 20000
$ws.Range("Q35").Value = "$/bandeja 15 kilos granel"
$ws.Range("R35").Value = "Región de O'Higgins"
$ws.Range("S35").Value = 1333
$ws.Range("T35").Value = 15
